# "Switching to Summer time" - refresh the Entsoe Consumption_Actual data:
# shift the timestamp base day from 45741 to 45744 and replace the
# consumption readings (column A) for existing rows, then append 8 new
# rows (49-56) continuing the same 15-minute cadence.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Actual Consumption (MW)" values for rows 2..56 (55 values)
$values = @(
    5576,5531,5542,5506,5466,5464,5451,5420,5384,5412,
    5410,5432,5413,5461,5453,5471,5561,5642,5652,5739,
    5877,5979,6169,6283,6549,6780,6879,6954,7072,7216,
    7277,7349,7358,7419,7371,7393,7340,7287,7305,7226,
    7139,7096,7103,7043,6970,6905,6863,6786,6654,6641,
    6623,6613,6585,6581,6510
)

$baseDay = 45744

# Reference number format for the timestamp column (already applied on B2:B48)
$dateFormat = $ws.Cells.Item(2, 2).NumberFormat

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2

    $ws.Cells.Item($row, 1).Value = $values[$i]

    $cellB = $ws.Cells.Item($row, 2)
    $cellB.Value = $baseDay + ($i / 96)
    $cellB.NumberFormat = $dateFormat
}
